$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string as TEXT (matches the source
# inlineStr cells) without leaving a lingering Text number-format on the
# cell -- format as Text just long enough to defeat Excels automatic
# "looks like a number" coercion, then clear the format back to General
# so the cells style index is unchanged from the original workbook.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$ws.Range('D2').Value = '44.293.82'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '2.239.27'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('E4').Value = '  +0.24%  '
Set-TextValue $ws.Range('D5') '307.13'
$ws.Range('E5').Value = '  -3.38%  '
Set-TextValue $ws.Range('D6') '93.47'
$ws.Range('E6').Value = '  -6.48%  '
Set-TextValue $ws.Range('D7') '0.571'
$ws.Range('E7').Value = '  -1.20%  '
Set-TextValue $ws.Range('D8') '1.00'
$ws.Range('E8').Value = '  +0.31%  '
Set-TextValue $ws.Range('D9') '0.522'
$ws.Range('E9').Value = '  -3.20%  '
Set-TextValue $ws.Range('D10') '34.20'
$ws.Range('E10').Value = '  -5.98%  '
Set-TextValue $ws.Range('D11') '0.0807'
$ws.Range('E11').Value = '  -2.57%  '
Set-TextValue $ws.Range('D12') '7.12'
$ws.Range('E12').Value = '  -4.43%  '
$ws.Range('E13').Value = '  -0.34%  '
$ws.Range('D14').Value = '2.337.39'
$ws.Range('E14').Value = '  +3.45%  '
Set-TextValue $ws.Range('D15') '0.827'
$ws.Range('E15').Value = '  -2.72%  '
Set-TextValue $ws.Range('D16') '13.41'
$ws.Range('E16').Value = '  -4.70%  '
$ws.Range('D17').Value = '44.012.33'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = '0.0₃0964'
$ws.Range('E18').Value = '  -2.02%  '
$ws.Range('B19').Value = 'InternetComputer(DFINITY)'
$ws.Range('C19').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D19') '12.04'
$ws.Range('E19').Value = '  -8.87%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D20') '6.31'
$ws.Range('E20').Value = '  -1.47%  '
Set-TextValue $ws.Range('D21') '65.63'
$ws.Range('E21').Value = '  -0.36%  '
$ws.Range('B22').Value = 'PancakeSwap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D22') '3.10'
$ws.Range('E22').Value = '  +2.75%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D23') '236.81'
$ws.Range('E23').Value = '  -1.46%  '
Set-TextValue $ws.Range('D24') '2.00'
$ws.Range('E24').Value = '  -2.76%  '
$ws.Range('E25').Value = '  -0.29%  '
Set-TextValue $ws.Range('D26') '39.49'
$ws.Range('E26').Value = '  +1.06%  '
$ws.Range('E27').Value = '  +3.62%  '
Set-TextValue $ws.Range('D28') '9.83'
$ws.Range('E28').Value = '  -4.13%  '
Set-TextValue $ws.Range('D29') '20.02'
$ws.Range('E29').Value = '  -0.93%  '
Set-TextValue $ws.Range('D30') '5.87'
$ws.Range('E30').Value = '  -3.32%  '
Set-TextValue $ws.Range('D31') '151.67'
$ws.Range('E31').Value = '  -2.31%  '
Set-TextValue $ws.Range('D32') '0.0794'
$ws.Range('E32').Value = '  -6.19%  '
$ws.Range('E33').Value = '  -2.78%  '
Set-TextValue $ws.Range('D34') '3.06'
$ws.Range('E34').Value = '  -12.12%  '
$ws.Range('E35').Value = '  +0.93%  '
$ws.Range('E36').Value = '  -2.80%  '
Set-TextValue $ws.Range('D37') '1.75'
$ws.Range('E37').Value = '  -9.26%  '
Set-TextValue $ws.Range('D38') '3.47'
$ws.Range('E38').Value = '  -1.96%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range('D39') '14.28'
$ws.Range('E39').Value = '  -8.16%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D40') '3.77'
$ws.Range('E40').Value = '  -4.97%  '
Set-TextValue $ws.Range('D41') '0.0297'
$ws.Range('E41').Value = '  -4.41%  '
$ws.Range('E42').Value = '  +0.30%  '
$ws.Range('D43').Value = '1.703.14'
$ws.Range('E43').Value = '  -2.05%  '
Set-TextValue $ws.Range('D44') '82.13'
$ws.Range('E44').Value = '  -3.27%  '
$ws.Range('E45').Value = '  -3.39%  '
Set-TextValue $ws.Range('D46') '99.28'
$ws.Range('E46').Value = '  -3.43%  '
Set-TextValue $ws.Range('D47') '4.91'
$ws.Range('E47').Value = '  -6.03%  '
$ws.Range('E48').Value = '  -2.87%  '
Set-TextValue $ws.Range('D49') '54.79'
$ws.Range('E49').Value = '  -4.18%  '
Set-TextValue $ws.Range('D50') '8.07'
$ws.Range('E50').Value = '  -2.62%  '
Set-TextValue $ws.Range('D51') '67.55'
$ws.Range('E51').Value = '  -5.46%  '
